$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("토요일")
$ws.Activate()

# Clear the reservation entry in row 2 (A2 numeric value, B2/C2 text values)
# but keep the existing number-format style already applied to B2/C2.
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()

# Update selection to A2:C2 with active cell C2
$ws.Range("A2:C2").Select()
$ws.Range("C2").Activate()
